$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.109.20'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.786.23'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.23'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.547'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '31.98'
$ws.Range('E8').Value = '  -2.72%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.292'
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0690'
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0943'
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = '2.043.03'
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.21'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('D14').Value = '1.789.93'
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').Value = '34.061.28'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '4.19'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '67.94'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '245.55'
$ws.Range('E19').Value = '  +0.18%  '
$ws.Range('D20').Value = '0.0₃0778'
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.05'
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '161.63'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('E26').Value = '  +0.32%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.32'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.60'
$ws.Range('E33').Value = '  +2.31%  '
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('D35').Value = '1.450.38'
$ws.Range('E35').Value = '  +3.37%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.644'
$ws.Range('E36').Value = '  -2.07%  '
$ws.Range('E37').Value = '  +1.90%  '
$ws.Range('E38').Value = '  +8.18%  '
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '80.03'
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.919'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.68'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.46'
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0508'
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').Value = '0.0₆0137'
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '107.52'
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('D50').Value = '1.944.60'
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('E51').Value = '  +0.27%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D49').Style = "Normal"
